$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("5001809-74.2017.8.21.0036", "9002336-55.2017.8.21.0036", "Migrado"),
    @("5002735-26.2018.8.21.0002", "9000271-58.2018.8.21.0002", "Migrado"),
    @("5002731-86.2018.8.21.0002", "9000433-53.2018.8.21.0002", "Migrado"),
    @("5035108-55.2014.8.21.0001", "0028693-44.2014.8.21.0001", "Digitalizado"),
    @("5035108-55.2014.8.21.0001", "0028693-44.2014.8.21.0001", "Digitalizado"),
    @("5035108-55.2014.8.21.0001", "0028693-44.2014.8.21.0001", "Digitalizado"),
    @("5035108-55.2014.8.21.0001", "0028693-44.2014.8.21.0001", "Digitalizado"),
    @("5009658-08.2018.8.21.0022", "9003683-34.2018.8.21.0022", "Migrado"),
    @("5002111-39.2019.8.21.0067", "9000482-59.2019.8.21.0067", "Migrado"),
    @("5002115-76.2019.8.21.0067", "9001255-07.2019.8.21.0067", "Migrado"),
    @("5001005-76.2018.8.21.0067", "9000856-12.2018.8.21.0067", "Migrado"),
    @("5002061-76.2020.8.21.0067", "9000182-63.2020.8.21.0067", "Migrado"),
    @("5002116-61.2019.8.21.0067", "9000433-18.2019.8.21.0067", "Migrado"),
    @("5002062-61.2020.8.21.0067", "9000400-91.2020.8.21.0067", "Migrado"),
    @("5001006-61.2018.8.21.0067", "9000817-15.2018.8.21.0067", "Migrado"),
    @("5002063-46.2020.8.21.0067", "9000178-26.2020.8.21.0067", "Migrado"),
    @("5002065-16.2020.8.21.0067", "9000388-77.2020.8.21.0067", "Migrado"),
    @("5001007-46.2018.8.21.0067", "9000790-32.2018.8.21.0067", "Migrado"),
    @("5002068-68.2020.8.21.0067", "9000458-94.2020.8.21.0067", "Migrado"),
    @("5000282-08.2011.8.21.0001", "0302109-66.2011.8.21.0001", "Digitalizado"),
    @("5034992-20.2012.8.21.0001", "0303805-06.2012.8.21.0001", "Digitalizado"),
    @("5009659-90.2018.8.21.0022", "9005551-47.2018.8.21.0022", "Migrado"),
    @("5002741-33.2018.8.21.0002", "9000456-96.2018.8.21.0002", "Migrado"),
    @("5008759-20.2011.8.21.0001", "0148964-87.2011.8.21.0001", "Digitalizado"),
    @("5019558-25.2011.8.21.0001", "0218063-47.2011.8.21.0001", "Digitalizado"),
    @("5125809-52.2020.8.21.0001", "9001529-57.2020.8.21.0027", "Migrado")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
